$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    $rng = $d.Paragraphs.Item($Index).Range
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        Write-Output "FAILED replace in paragraph $Index : $OldText"
    }
}

# --- Paragraph 7 : "Para este punto se sobreentiende..." ---
Replace-InParagraph 7 "Rick y Morty tienen" "Morty tienen"
Replace-InParagraph 7 "jugador (Rick, Morty o los perros)" "jugador (Morty o los perros)"
Replace-InParagraph 7 "saltando (como el perro no tendrá esta capacidad, simplemente se va a chocar y se eliminará) de esta forma" "bordeándolo y como hay un margen de distancia con el perro, este no tendrá tiempo de esquivarlo y se eliminará, de esta forma"

# --- Paragraph 8 : "Si el jugador se deja tocar..." ---
Replace-InParagraph 8 "perros, esto le restará" "perros, este le restará"

# --- Paragraph 9 : "Este nivel durará..." ---
Replace-InParagraph 9 "Este nivel durará aproximadamente 1 minuto, y para pasar al siguiente nivel" "Para pasar al siguiente nivel"

# --- Paragraph 10 : "Si cuando termine el tiempo indicado..." (full replace) ---
Replace-InParagraph 10 "Si cuando termine el tiempo indicado, el jugador no acumuló la cantidad mínima de puntos necesarios para pasar al siguiente nivel, morirá y le toca iniciar desde el principio." "Si acumuló los puntos necesarios, no solo aparecerá un portal sino varios y el jugador tendrán que encontrar el verdadero, el que lo va a llevar al siguiente nivel. Van a ser 3 portales, uno que no hará nada, el verdadero y el que lo va a llevar al inicio del juego."

# --- Paragraph 11 : "En este nuevo nivel..." (full replace) ---
Replace-InParagraph 11 "En este nuevo nivel el jugador se encontrará con que, si bien la cantidad de perros va a estar disminuida, estos ya van a tener armas para atacarlo (no solo correr detrás de él), además habrá paredes que van a servir como escudo al jugador para protegerse de las balas ya que, si estas lo tocan, le restará vida. Pero estas paredes a la vez también tendrán vida y una vez esta se agote, se desaparecerá." "En este nuevo nivel solo habrá un perro que se va a mover rápido tratando de tocar al jugador y en ciertos momentos saltará sobre él para tratar de aplastarlo, el jugador tiene que esquivarlo y va a tener una pistola para atacar al perro grande (ambos se moverán de forma horizontal)."

# --- Paragraph 12 : "Cabe aclarar..." (full replace) ---
Replace-InParagraph 12 "Cabe aclarar que la cantidad de vida que tenga el jugador no se reinicia entre niveles, por lo que le tocará cuidarla muy bien. " "Hay que aclarar que la vida del jugador no se reinicia cuando pase de nivel, por lo que le tocará cuidarla muy bien."

# --- Remove the trailing paragraphs 13-17 (old "Para el jugador defenderse", "Para pasar al tercer...",
#     "Es posible...", and two blank paragraphs) so the document ends right after paragraph 12. ---
for ($n = 1; $n -le 4; $n++) {
    $d.Paragraphs.Item(13).Range.Delete()
}
# The very last paragraph mark in a document can't be removed by deleting its own Range
# (Word always keeps a final paragraph mark), so merge it into the prior paragraph instead:
# collapse to the end of paragraph 12 and extend one character to grab the trailing pilcrow.
$tail = $d.Paragraphs.Item(12).Range
$tail.Collapse(0)
$null = $tail.MoveEnd(1, 1)
$tail.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
